# Update "想去人数" (F column) values across all sheets to match
# the regenerated gh-pages output at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 3554
$ws.Cells.Item(5, 6).Value = 8316
$ws.Cells.Item(7, 6).Value = 122
$ws.Cells.Item(8, 6).Value = 2224
$ws.Cells.Item(10, 6).Value = 94
$ws.Cells.Item(12, 6).Value = 643
$ws.Cells.Item(14, 6).Value = 7379
$ws.Cells.Item(16, 6).Value = 7641
$ws.Cells.Item(18, 6).Value = 57627
$ws.Cells.Item(19, 6).Value = 57627
$ws.Cells.Item(20, 6).Value = 4748
$ws.Cells.Item(22, 6).Value = 940
$ws.Cells.Item(24, 6).Value = 108
$ws.Cells.Item(25, 6).Value = 925
$ws.Cells.Item(28, 6).Value = 5293
$ws.Cells.Item(29, 6).Value = 600
$ws.Cells.Item(30, 6).Value = 111
$ws.Cells.Item(31, 6).Value = 47
$ws.Cells.Item(32, 6).Value = 905
$ws.Cells.Item(33, 6).Value = 1365
$ws.Cells.Item(34, 6).Value = 1894
$ws.Cells.Item(36, 6).Value = 182
$ws.Cells.Item(37, 6).Value = 226
$ws.Cells.Item(38, 6).Value = 1086
$ws.Cells.Item(40, 6).Value = 728
$ws.Cells.Item(42, 6).Value = 780
$ws.Cells.Item(43, 6).Value = 260
$ws.Cells.Item(44, 6).Value = 173
$ws.Cells.Item(46, 6).Value = 17
$ws.Cells.Item(47, 6).Value = 200
$ws.Cells.Item(49, 6).Value = 58
$ws.Cells.Item(50, 6).Value = 2485

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 200
$ws.Cells.Item(5, 6).Value = 63
$ws.Cells.Item(9, 6).Value = 181
$ws.Cells.Item(10, 6).Value = 7614
$ws.Cells.Item(11, 6).Value = 125
$ws.Cells.Item(14, 6).Value = 2
$ws.Cells.Item(26, 6).Value = 2
$ws.Cells.Item(30, 6).Value = 3
$ws.Cells.Item(35, 6).Value = 3

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 6).Value = 79
$ws.Cells.Item(4, 6).Value = 2363
$ws.Cells.Item(5, 6).Value = 1595
$ws.Cells.Item(7, 6).Value = 680
$ws.Cells.Item(8, 6).Value = 2412
$ws.Cells.Item(10, 6).Value = 1759
$ws.Cells.Item(12, 6).Value = 112
$ws.Cells.Item(15, 6).Value = 267
$ws.Cells.Item(16, 6).Value = 2335
$ws.Cells.Item(17, 6).Value = 61
$ws.Cells.Item(18, 6).Value = 504

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 3554
$ws.Cells.Item(3, 6).Value = 2363
$ws.Cells.Item(5, 6).Value = 8316
$ws.Cells.Item(6, 6).Value = 680
$ws.Cells.Item(8, 6).Value = 122
$ws.Cells.Item(9, 6).Value = 267
$ws.Cells.Item(11, 6).Value = 643
$ws.Cells.Item(13, 6).Value = 7379
$ws.Cells.Item(14, 6).Value = 7641
$ws.Cells.Item(15, 6).Value = 57627
$ws.Cells.Item(16, 6).Value = 200
$ws.Cells.Item(18, 6).Value = 4748
$ws.Cells.Item(20, 6).Value = 940
$ws.Cells.Item(24, 6).Value = 5293
$ws.Cells.Item(25, 6).Value = 600
$ws.Cells.Item(26, 6).Value = 111
$ws.Cells.Item(27, 6).Value = 905
$ws.Cells.Item(28, 6).Value = 1365
$ws.Cells.Item(29, 6).Value = 125
$ws.Cells.Item(30, 6).Value = 504
$ws.Cells.Item(33, 6).Value = 182
$ws.Cells.Item(35, 6).Value = 728
$ws.Cells.Item(36, 6).Value = 780
$ws.Cells.Item(37, 6).Value = 260
$ws.Cells.Item(40, 6).Value = 2
$ws.Cells.Item(45, 6).Value = 200
$ws.Cells.Item(48, 6).Value = 43
$ws.Cells.Item(49, 6).Value = 2485

